# Merge Look-up tables - new BC Transformer and workspace
#
# This lookup table's source sheet was renamed from "Keyword" to "Country"
# as part of merging it into the new BC Transformer workspace, and the
# active selection left on the sheet moved to A19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Country"
$ws.Range("A19").Select()
